# Update the "New M1" percent-change column (F) on the "Sheet" worksheet.
# The percentage applied to M1 (column C) changes from 70% to 54%,
# "depending on users input" per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$newValues = @{
    2  = 24.3
    3  = 46.44
    4  = 25.92
    5  = 30.78
    6  = 45.90000000000001
    7  = 40.5
    8  = 31.32
    9  = 18.9
    10 = 21.6
    11 = 49.68000000000001
    12 = 51.3
    13 = 48.06
    14 = 28.08
    15 = 31.32
    16 = 33.48
    17 = 53.46
    18 = 23.22
    19 = 38.34
    20 = 45.36
    21 = 34.02
    22 = 36.72
    23 = 27
    24 = 24.3
    25 = 29.16
    26 = 30.24
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
